$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 70.2
$ws.Range("I5").Value = 9.5
$ws.Range("J5").Value = 110.666664
$ws.Range("K5").Value = 9.5
$ws.Range("L5").Value = 110.666664
$ws.Range("M5").Value = 105.5
$ws.Range("N5").Value = -340.666664
$ws.Range("H40").Value = 1361.5385
$ws.Range("I40").Value = 1375
$ws.Range("J40").Value = 1200
$ws.Range("K40").Value = 1375
$ws.Range("L40").Value = 1200
$ws.Range("M40").Value = -1200
$ws.Range("N40").Value = -1550
$ws.Range("H70").Value = 5188.375
$ws.Range("I70").Value = 1400
$ws.Range("J70").Value = 6451.1665
$ws.Range("K70").Value = 4200
$ws.Range("L70").Value = 19353.4995
$ws.Range("M70").Value = -3930
$ws.Range("N70").Value = -19893.4995
$ws.Range("H73").Value = 5188.375
$ws.Range("I73").Value = 1400
$ws.Range("J73").Value = 6451.1665
$ws.Range("K73").Value = 4200
$ws.Range("L73").Value = 19353.4995
$ws.Range("M73").Value = -3264
$ws.Range("N73").Value = -21225.4995
$ws.Range("H132").Value = 3704904.2
$ws.Range("I132").Value = 1285.3462
$ws.Range("J132").Value = 27778428
$ws.Range("K132").Value = 3856.0386
$ws.Range("L132").Value = 83335284
$ws.Range("M132").Value = -1326.0386
$ws.Range("N132").Value = -83340344

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 13904.852
$ws.Range("I2").Value = 18395.05
$ws.Range("J2").Value = 1075.7142
$ws.Range("K2").Value = 18395.05
$ws.Range("L2").Value = 1075.7142
$ws.Range("M2").Value = -18282.05
$ws.Range("N2").Value = -1301.7142
$ws.Range("H88").Value = 3923.05
$ws.Range("I88").Value = 1937.1
$ws.Range("J88").Value = 5909
$ws.Range("K88").Value = 1937.1
$ws.Range("L88").Value = 5909
$ws.Range("M88").Value = -1531.1
$ws.Range("N88").Value = -6721
$ws.Range("H91").Value = 3923.05
$ws.Range("I91").Value = 1937.1
$ws.Range("J91").Value = 5909
$ws.Range("K91").Value = 1937.1
$ws.Range("L91").Value = 5909
$ws.Range("M91").Value = -533.0999999999999
$ws.Range("N91").Value = -8717
$ws.Range("H116").Value = 13904.852
$ws.Range("I116").Value = 18395.05
$ws.Range("J116").Value = 1075.7142
$ws.Range("K116").Value = 18395.05
$ws.Range("L116").Value = 1075.7142
$ws.Range("M116").Value = -16101.05
$ws.Range("N116").Value = -5663.7142

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 13904.852
$ws.Range("I3").Value = 18395.05
$ws.Range("J3").Value = 1075.7142
$ws.Range("K3").Value = 18395.05
$ws.Range("L3").Value = 1075.7142
$ws.Range("M3").Value = -18281.05
$ws.Range("N3").Value = -1303.7142
$ws.Range("H86").Value = 1940.2
$ws.Range("I86").Value = 1946.0103
$ws.Range("J86").Value = 1752.3334
$ws.Range("K86").Value = 1946.0103
$ws.Range("L86").Value = 1752.3334
$ws.Range("M86").Value = -823.0102999999999
$ws.Range("N86").Value = -3998.3334
$ws.Range("H89").Value = 1940.2
$ws.Range("I89").Value = 1946.0103
$ws.Range("J89").Value = 1752.3334
$ws.Range("K89").Value = 9730.0515
$ws.Range("L89").Value = 8761.666999999999
$ws.Range("M89").Value = -4114.0515
$ws.Range("N89").Value = -19993.667
$ws.Range("H105").Value = 2171.6667
$ws.Range("I105").Value = 2087.7778
$ws.Range("J105").Value = 2255.5557
$ws.Range("K105").Value = 2087.7778
$ws.Range("L105").Value = 2255.5557
$ws.Range("M105").Value = -340.7777999999998
$ws.Range("N105").Value = -5749.5557
$ws.Range("H106").Value = 47500
$ws.Range("J106").Value = 47500
$ws.Range("L106").Value = 47500
$ws.Range("N106").Value = -50024

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 10434.375
$ws.Range("I86").Value = 18464.428
$ws.Range("J86").Value = 4188.778
$ws.Range("K86").Value = 18464.428
$ws.Range("L86").Value = 4188.778
$ws.Range("M86").Value = -17341.428
$ws.Range("N86").Value = -6434.778
$ws.Range("H89").Value = 10434.375
$ws.Range("I89").Value = 18464.428
$ws.Range("J89").Value = 4188.778
$ws.Range("K89").Value = 92322.14
$ws.Range("L89").Value = 20943.89
$ws.Range("M89").Value = -86706.14
$ws.Range("N89").Value = -32175.89
$ws.Range("H133").Value = 13263.158
$ws.Range("J133").Value = 13263.158
$ws.Range("L133").Value = 13263.158
$ws.Range("N133").Value = -18323.158

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1616.6666
$ws.Range("I34").Value = 583.3333
$ws.Range("J34").Value = 2133.3333
$ws.Range("K34").Value = 1749.9999
$ws.Range("L34").Value = 6399.999899999999
$ws.Range("M34").Value = -1665.9999
$ws.Range("N34").Value = -6567.999899999999
$ws.Range("H39").Value = 2316.1428
$ws.Range("J39").Value = 3066.6667
$ws.Range("L39").Value = 9200.000100000001
$ws.Range("N39").Value = -9788.000100000001
$ws.Range("H55").Value = 2750.3333
$ws.Range("I55").Value = 1004
$ws.Range("J55").Value = 2909.0908
$ws.Range("K55").Value = 3012
$ws.Range("L55").Value = 8727.2724
$ws.Range("M55").Value = -2835
$ws.Range("N55").Value = -9081.2724

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 10867.889
$ws.Range("I80").Value = 8145.6
$ws.Range("J80").Value = 11486.591
$ws.Range("K80").Value = 8145.6
$ws.Range("L80").Value = 11486.591
$ws.Range("M80").Value = -7147.6
$ws.Range("N80").Value = -13482.591
$ws.Range("H83").Value = 10867.889
$ws.Range("I83").Value = 8145.6
$ws.Range("J83").Value = 11486.591
$ws.Range("K83").Value = 40728
$ws.Range("L83").Value = 57432.955
$ws.Range("M83").Value = -35736
$ws.Range("N83").Value = -67416.955

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H62").Value = 756000
$ws.Range("I62").Value = 24000
$ws.Range("J62").Value = 1000000
$ws.Range("K62").Value = 24000
$ws.Range("L62").Value = 1000000
$ws.Range("M62").Value = -23376
$ws.Range("N62").Value = -1001248
$ws.Range("H63").Value = 1000000
$ws.Range("J63").Value = 1000000
$ws.Range("L63").Value = 1000000
$ws.Range("N63").Value = -1001498
$ws.Range("H65").Value = 756000
$ws.Range("I65").Value = 24000
$ws.Range("J65").Value = 1000000
$ws.Range("K65").Value = 72000
$ws.Range("L65").Value = 3000000
$ws.Range("M65").Value = -68880
$ws.Range("N65").Value = -3006240
$ws.Range("H66").Value = 1000000
$ws.Range("J66").Value = 1000000
$ws.Range("L66").Value = 3000000
$ws.Range("N66").Value = -3007488
$ws.Range("H68").Value = 144588.58
$ws.Range("I68").Value = 2066.6667
$ws.Range("J68").Value = 251480
$ws.Range("K68").Value = 2066.6667
$ws.Range("L68").Value = 251480
$ws.Range("M68").Value = -1317.6667
$ws.Range("N68").Value = -252978
$ws.Range("H71").Value = 144588.58
$ws.Range("I71").Value = 2066.6667
$ws.Range("J71").Value = 251480
$ws.Range("K71").Value = 10333.3335
$ws.Range("L71").Value = 1257400
$ws.Range("M71").Value = -6589.333500000001
$ws.Range("N71").Value = -1264888

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 3500
$ws.Range("I39").Value = 3500
$ws.Range("K39").Value = 3500
$ws.Range("M39").Value = -3087
$ws.Range("H43").Value = 5794
$ws.Range("I43").Value = 5000
$ws.Range("J43").Value = 5992.5
$ws.Range("K43").Value = 5000
$ws.Range("L43").Value = 5992.5
$ws.Range("M43").Value = -4851
$ws.Range("N43").Value = -6290.5
$ws.Range("H136").Value = 19076.846
$ws.Range("I136").Value = 12798.4
$ws.Range("J136").Value = 40005
$ws.Range("K136").Value = 38395.2
$ws.Range("L136").Value = 120015
$ws.Range("M136").Value = -35845.2
$ws.Range("N136").Value = -125115
